# Updated cryptos list - refresh Price (D) and Volume(1h) (E) figures,
# and fix the swapped Polkadot/WrappedEther and InjectiveProtocol/Kaspa rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.425.44"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.912.58"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.661"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.27%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.62"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.353"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.39%  "
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0997"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.30%  "
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.909.00"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.89"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").Value = "35.393.63"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").Value = "0.0₃0824"
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.51"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.85"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +21.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +6.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0567"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.939"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.92%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.75"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0209"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0652"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +17.15%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.37"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "90.39"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "1.345.50"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.13"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +38.19%  "
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "2.094.29"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0699"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.79%  "
